$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 13579.444
$ws.Range("I58").Value = 725
$ws.Range("J58").Value = 15186.25
$ws.Range("K58").Value = 2175
$ws.Range("L58").Value = 45558.75
$ws.Range("M58").Value = -2025
$ws.Range("N58").Value = -45858.75
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
# Row 86
$ws.Range("H86").Value = 3800.7778
$ws.Range("I86").Value = 2267.6667
$ws.Range("J86").Value = 4567.3335
$ws.Range("K86").Value = 2267.6667
$ws.Range("L86").Value = 4567.3335
$ws.Range("M86").Value = -1144.6667
$ws.Range("N86").Value = -6813.3335
# Row 89
$ws.Range("H89").Value = 3800.7778
$ws.Range("I89").Value = 2267.6667
$ws.Range("J89").Value = 4567.3335
$ws.Range("K89").Value = 11338.3335
$ws.Range("L89").Value = 22836.6675
$ws.Range("M89").Value = -5722.333500000001
$ws.Range("N89").Value = -34068.6675
# Row 98
$ws.Range("H98").Value = 4562.154
$ws.Range("I98").Value = 2430.8
$ws.Range("J98").Value = 11666.667
$ws.Range("K98").Value = 2430.8
$ws.Range("L98").Value = 11666.667
$ws.Range("M98").Value = -932.8000000000002
$ws.Range("N98").Value = -14662.667
# Row 103
$ws.Range("H103").Value = 14656.667
$ws.Range("I103").Value = 588
$ws.Range("J103").Value = 85000
$ws.Range("K103").Value = 1764
$ws.Range("L103").Value = 255000
$ws.Range("M103").Value = -1178
$ws.Range("N103").Value = -256172
# Row 112
$ws.Range("H112").Value = 2145.4443
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 2394.5334
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 7183.600199999999
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -9399.600199999999
# Row 122
$ws.Range("H122").Value = 4562.154
$ws.Range("I122").Value = 2430.8
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 7292.400000000001
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -4842.400000000001
$ws.Range("N122").Value = -39900.001
# Row 129
$ws.Range("H129").Value = 1016.7059
$ws.Range("I129").Value = 499
$ws.Range("J129").Value = 1032.3939
$ws.Range("K129").Value = 1497
$ws.Range("L129").Value = 3097.1817
$ws.Range("M129").Value = 3503
$ws.Range("N129").Value = -13097.1817
# Row 132
$ws.Range("H132").Value = 30426976
$ws.Range("I132").Value = 35859156
$ws.Range("J132").Value = 6771.2
$ws.Range("K132").Value = 107577468
$ws.Range("L132").Value = 20313.6
$ws.Range("M132").Value = -107574938
$ws.Range("N132").Value = -25373.6
# Row 138
$ws.Range("H138").Value = 3936.4023
$ws.Range("I138").Value = 2852.1428
$ws.Range("J138").Value = 4144.3423
$ws.Range("K138").Value = 8556.428400000001
$ws.Range("L138").Value = 12433.0269
$ws.Range("M138").Value = -3416.428400000001
$ws.Range("N138").Value = -22713.0269

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2487.05
$ws.Range("I61").Value = 1641.5834
$ws.Range("J61").Value = 3755.25
$ws.Range("K61").Value = 1641.5834
$ws.Range("L61").Value = 3755.25
$ws.Range("M61").Value = -1429.5834
$ws.Range("N61").Value = -4179.25
# Row 74
$ws.Range("H74").Value = 1648.7441
$ws.Range("I74").Value = 1244.9117
$ws.Range("J74").Value = 3174.3333
$ws.Range("K74").Value = 1244.9117
$ws.Range("L74").Value = 3174.3333
$ws.Range("M74").Value = -370.9117000000001
$ws.Range("N74").Value = -4922.3333
# Row 77
$ws.Range("H77").Value = 1648.7441
$ws.Range("I77").Value = 1244.9117
$ws.Range("J77").Value = 3174.3333
$ws.Range("K77").Value = 6224.558500000001
$ws.Range("L77").Value = 15871.6665
$ws.Range("M77").Value = -1856.558500000001
$ws.Range("N77").Value = -24607.6665
# Row 136
$ws.Range("H136").Value = 2487.05
$ws.Range("I136").Value = 1641.5834
$ws.Range("J136").Value = 3755.25
$ws.Range("K136").Value = 4924.7502
$ws.Range("L136").Value = 11265.75
$ws.Range("M136").Value = -2374.7502
$ws.Range("N136").Value = -16365.75
# Row 139
$ws.Range("H139").Value = 41527.668
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 41527.668
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 41527.668
$ws.Range("N139").Value = -51807.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2667.8462
$ws.Range("I20").Value = 2816.1765
$ws.Range("J20").Value = 2387.6667
$ws.Range("K20").Value = 2816.1765
$ws.Range("L20").Value = 2387.6667
$ws.Range("M20").Value = -2569.1765
$ws.Range("N20").Value = -2881.6667
# Row 21
$ws.Range("H21").Value = 22156.166
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 22156.166
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 22156.166
$ws.Range("N21").Value = -22628.166
# Row 64
$ws.Range("H64").Value = 207.9
$ws.Range("I64").Value = 172.4
$ws.Range("J64").Value = 243.4
$ws.Range("K64").Value = 172.4
$ws.Range("L64").Value = 243.4
$ws.Range("M64").Value = 52.59999999999999
$ws.Range("N64").Value = -693.4
# Row 67
$ws.Range("H67").Value = 207.9
$ws.Range("I67").Value = 172.4
$ws.Range("J67").Value = 243.4
$ws.Range("K67").Value = 172.4
$ws.Range("L67").Value = 243.4
$ws.Range("M67").Value = 607.6
$ws.Range("N67").Value = -1803.4
# Row 99
$ws.Range("H99").Value = 3265.7693
$ws.Range("I99").Value = 1136.6666
$ws.Range("J99").Value = 3543.4783
$ws.Range("K99").Value = 1136.6666
$ws.Range("L99").Value = 3543.4783
$ws.Range("M99").Value = 361.3334
$ws.Range("N99").Value = -6539.478300000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 7107.8335
$ws.Range("I99").Value = 5215
$ws.Range("J99").Value = 9000.666999999999
$ws.Range("K99").Value = 5215
$ws.Range("L99").Value = 9000.666999999999
$ws.Range("M99").Value = -3717
$ws.Range("N99").Value = -11996.667
# Row 126
$ws.Range("H126").Value = 7107.8335
$ws.Range("I126").Value = 5215
$ws.Range("J126").Value = 9000.666999999999
$ws.Range("K126").Value = 15645
$ws.Range("L126").Value = 27002.001
$ws.Range("M126").Value = -13175
$ws.Range("N126").Value = -31942.001
# Row 132
$ws.Range("H132").Value = 2837.1191
$ws.Range("I132").Value = 2425.9656
$ws.Range("J132").Value = 3754.3076
$ws.Range("K132").Value = 7277.8968
$ws.Range("L132").Value = 11262.9228
$ws.Range("M132").Value = -4747.8968
$ws.Range("N132").Value = -16322.9228
# Row 134
$ws.Range("H134").Value = 4313.7954
$ws.Range("I134").Value = 4715.1562
$ws.Range("J134").Value = 3243.5
$ws.Range("K134").Value = 14145.4686
$ws.Range("L134").Value = 9730.5
$ws.Range("M134").Value = -11610.4686
$ws.Range("N134").Value = -14800.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 48.166668
$ws.Range("I38").Value = 20.90909
$ws.Range("J38").Value = 91
$ws.Range("K38").Value = 62.72727
$ws.Range("L38").Value = 273
$ws.Range("M38").Value = 284.27273
$ws.Range("N38").Value = -967
# Row 39
$ws.Range("H39").Value = 12594.823
$ws.Range("I39").Value = 16000
$ws.Range("J39").Value = 12382
$ws.Range("K39").Value = 48000
$ws.Range("L39").Value = 37146
$ws.Range("M39").Value = -47706
$ws.Range("N39").Value = -37734
# Row 107
$ws.Range("H107").Value = 1789.3636
$ws.Range("I107").Value = 440
$ws.Range("J107").Value = 2913.8333
$ws.Range("K107").Value = 1320
$ws.Range("L107").Value = 8741.499899999999
$ws.Range("M107").Value = 600
$ws.Range("N107").Value = -12581.4999
# Row 122
$ws.Range("H122").Value = 2812.2363
$ws.Range("I122").Value = 630.9286
$ws.Range("J122").Value = 3557.0732
$ws.Range("K122").Value = 5678.3574
$ws.Range("L122").Value = 32013.6588
$ws.Range("M122").Value = -3228.3574
$ws.Range("N122").Value = -36913.6588
# Row 132
$ws.Range("H132").Value = 3307.7273
$ws.Range("I132").Value = 1024.5
$ws.Range("J132").Value = 4163.9375
$ws.Range("K132").Value = 9220.5
$ws.Range("L132").Value = 37475.4375
$ws.Range("M132").Value = -6690.5
$ws.Range("N132").Value = -42535.4375
# Row 137
$ws.Range("H137").Value = 2782.15
$ws.Range("I137").Value = 2734
$ws.Range("J137").Value = 2926.6
$ws.Range("K137").Value = 8202
$ws.Range("L137").Value = 8779.799999999999
$ws.Range("M137").Value = -3102
$ws.Range("N137").Value = -18979.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3956.4375
$ws.Range("I122").Value = 3002.818
$ws.Range("J122").Value = 6054.4
$ws.Range("K122").Value = 9008.454000000002
$ws.Range("L122").Value = 18163.2
$ws.Range("M122").Value = -6558.454000000002
$ws.Range("N122").Value = -23063.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2312.5
$ws.Range("I46").Value = 2420
$ws.Range("J46").Value = 2263.6365
$ws.Range("K46").Value = 2420
$ws.Range("L46").Value = 2263.6365
$ws.Range("M46").Value = -2232
$ws.Range("N46").Value = -2639.6365
# Row 61
$ws.Range("H61").Value = 1512.0714
$ws.Range("I61").Value = 1438.6666
$ws.Range("J61").Value = 1952.5
$ws.Range("K61").Value = 1438.6666
$ws.Range("L61").Value = 1952.5
$ws.Range("M61").Value = -1236.6666
$ws.Range("N61").Value = -2356.5
# Row 100
$ws.Range("H100").Value = 2470
$ws.Range("I100").Value = 2258
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2258
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1717
$ws.Range("N100").Value = -4082
# Row 113
$ws.Range("H113").Value = 1512.0714
$ws.Range("I113").Value = 1438.6666
$ws.Range("J113").Value = 1952.5
$ws.Range("K113").Value = 1438.6666
$ws.Range("L113").Value = 1952.5
$ws.Range("M113").Value = 731.3334
$ws.Range("N113").Value = -6292.5
# Row 122
$ws.Range("H122").Value = 5624.9165
$ws.Range("I122").Value = 3099.8333
$ws.Range("J122").Value = 8150
$ws.Range("K122").Value = 9299.499899999999
$ws.Range("L122").Value = 24450
$ws.Range("M122").Value = -6849.499899999999
$ws.Range("N122").Value = -29350

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 10405
$ws.Range("I113").Value = 20340.4
$ws.Range("J113").Value = 469.6
$ws.Range("K113").Value = 61021.2
$ws.Range("L113").Value = 1408.8
$ws.Range("M113").Value = -58851.2
$ws.Range("N113").Value = -5748.8
